$wb = $excel.ActiveWorkbook

# FEB-22 is the active sheet (activeTab = 1, second sheet)
$ws = $wb.Worksheets.Item("FEB-22")

# Update the existing comment text at D21 to insert "and " before "waiting"
$ws.Range("D21").Value = "1. Customization due to Multi-Factor Authentication  has been implemented at DAILY-tasks (DRS, WARRANTY, SAW_DISCOUNT, SERVICE_ORDER_MANAGEMENT BULK, COMMUNICATION BOARD CLEAR TASK), and waiting for the OTP process"

# Fill in the new daily-track entry on row 22
$ws.Range("A22").Value = 15
$ws.Range("B22").Value = 44609
$ws.Range("C22").Value = "RPA RLOGIC"
$ws.Range("D22").Value = "1. Research and Development is going on general Captcha issue"
$ws.Range("E22").Value = 0.1
$ws.Range("F22").Value = "WIP"

# Update the selected cell on the active sheet
$ws.Range("D28").Select()
